$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 36 (record #34) is currently blank, formatted the same way as the other
# "支出"/"生活费" expense rows above it (e.g. row 33). Copy that row's
# formatting + values down, then only change what's different for this entry:
# the date and the note/remark text. (D36 amount of 400 is already correct
# after the copy, since it matches row 33's amount.)
$ws.Range("C33:G33").Copy($ws.Range("C36:G36"))

$ws.Range("E36").Value = (Get-Date -Year 2018 -Month 4 -Day 20 -Hour 0 -Minute 0 -Second 0)
$ws.Range("G36").Value = "生活费(4/21-4/30)"
# Re-assert F36's value (unchanged text, "生活费") so the SUMIFS totals that
# key off column F ("J9") see the newly-populated row and recalc correctly.
$ws.Range("F36").Value = "生活费"

# Update the view state: scroll position and active selection, matching the diff
$excel.ActiveWindow.ScrollRow = 28
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G37").Select()
